$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 32
$ws.Range("H32").Value = 14998.75
$ws.Range("J32").Value = 13331.667
$ws.Range("L32").Value = 13331.667
$ws.Range("N32").Value = -13983.667

# Row 51
$ws.Range("H51").Value = 9671.233
$ws.Range("J51").Value = 10264.708
$ws.Range("L51").Value = 10264.708
$ws.Range("N51").Value = -11232.708

# Row 70
$ws.Range("H70").Value = 23900
$ws.Range("J70").Value = 35250
$ws.Range("L70").Value = 105750
$ws.Range("N70").Value = -106290

# Row 73
$ws.Range("H73").Value = 23900
$ws.Range("J73").Value = 35250
$ws.Range("L73").Value = 105750
$ws.Range("N73").Value = -107622

# Row 92
$ws.Range("H92").Value = 718.1739
$ws.Range("I92").Value = 423.27777
$ws.Range("K92").Value = 423.27777
$ws.Range("M92").Value = 824.7222300000001

# Row 138
$ws.Range("H138").Value = 4736.391
$ws.Range("I138").Value = 1373.1765
$ws.Range("K138").Value = 4119.529500000001
$ws.Range("M138").Value = 1020.470499999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 63
$ws.Range("H63").Value = 4131.8335
$ws.Range("I63").Value = 3958.2
$ws.Range("K63").Value = 3958.2
$ws.Range("M63").Value = -3272.2

# Row 66
$ws.Range("H66").Value = 4131.8335
$ws.Range("I66").Value = 3958.2
$ws.Range("K66").Value = 19791
$ws.Range("M66").Value = -16359

# Row 74
$ws.Range("I74").Value = 100010570
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 100010570
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -100009696
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("I77").Value = 100010570
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 500052850
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -500048482
$ws.Range("N77").ClearContents()

# Row 97
$ws.Range("H97").Value = 571.3889
$ws.Range("I97").Value = 560.1429000000001
$ws.Range("K97").Value = 560.1429000000001
$ws.Range("M97").Value = -64.14290000000005

# Row 132
$ws.Range("H132").Value = 2442724.5
$ws.Range("I132").Value = 3451352
$ws.Range("K132").Value = 10354056
$ws.Range("M132").Value = -10351526

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 62
$ws.Range("H62").Value = 95063.336
$ws.Range("I62").Value = 81900
$ws.Range("K62").Value = 81900
$ws.Range("M62").Value = -81214

# Row 65
$ws.Range("H65").Value = 95063.336
$ws.Range("I65").Value = 81900
$ws.Range("K65").Value = 245700
$ws.Range("M65").Value = -242268

# Row 86
$ws.Range("H86").Value = 2004.2307
$ws.Range("I86").Value = 2009.5555
$ws.Range("K86").Value = 2009.5555
$ws.Range("M86").Value = -886.5554999999999

# Row 89
$ws.Range("H89").Value = 2004.2307
$ws.Range("I89").Value = 2009.5555
$ws.Range("K89").Value = 10047.7775
$ws.Range("M89").Value = -4431.7775

# Row 105
$ws.Range("H105").Value = 2589.8333
$ws.Range("I105").Value = 2507.8
$ws.Range("K105").Value = 2507.8
$ws.Range("M105").Value = -760.8000000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 12291.963
$ws.Range("I31").Value = 8002.75
$ws.Range("J31").Value = 15723.333
$ws.Range("K31").Value = 8002.75
$ws.Range("L31").Value = 15723.333
$ws.Range("M31").Value = -7707.75
$ws.Range("N31").Value = -16313.333

# Row 34
$ws.Range("H34").Value = 12291.963
$ws.Range("I34").Value = 8002.75
$ws.Range("J34").Value = 15723.333
$ws.Range("K34").Value = 8002.75
$ws.Range("L34").Value = 15723.333
$ws.Range("M34").Value = -7800.75
$ws.Range("N34").Value = -16127.333

# Row 58
$ws.Range("H58").Value = 33341528
$ws.Range("J58").Value = 2465
$ws.Range("L58").Value = 2465
$ws.Range("N58").Value = -2871

# Row 99
$ws.Range("H99").Value = 2260.2856
$ws.Range("I99").Value = 2150.5
$ws.Range("J99").Value = 2304.2
$ws.Range("K99").Value = 2150.5
$ws.Range("L99").Value = 2304.2
$ws.Range("M99").Value = -652.5
$ws.Range("N99").Value = -5300.2

# Row 122
$ws.Range("H122").Value = 5351.2666
$ws.Range("I122").Value = 5269.2856
$ws.Range("K122").Value = 15807.8568
$ws.Range("M122").Value = -13357.8568

# Row 126
$ws.Range("H126").Value = 2260.2856
$ws.Range("I126").Value = 2150.5
$ws.Range("J126").Value = 2304.2
$ws.Range("K126").Value = 6451.5
$ws.Range("L126").Value = 6912.599999999999
$ws.Range("M126").Value = -3981.5
$ws.Range("N126").Value = -11852.6

# Row 132
$ws.Range("H132").Value = 200003740
$ws.Range("I132").Value = 200003740
$ws.Range("K132").Value = 600011220
$ws.Range("M132").Value = -600008690

# Row 136
$ws.Range("H136").Value = 33341528
$ws.Range("J136").Value = 2465
$ws.Range("L136").Value = 7395
$ws.Range("N136").Value = -12495

# Row 141
$ws.Range("H141").Value = 304013.44
$ws.Range("J141").Value = 304013.44
$ws.Range("L141").Value = 304013.44
$ws.Range("N141").Value = -314373.44

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 86
$ws.Range("H86").Value = 805.4545000000001
$ws.Range("I86").Value = 659.8
$ws.Range("K86").Value = 1979.4
$ws.Range("M86").Value = -793.3999999999999

# Row 89
$ws.Range("H89").Value = 805.4545000000001
$ws.Range("I89").Value = 659.8
$ws.Range("K89").Value = 5938.2
$ws.Range("M89").Value = -10.19999999999982

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 3673.625
$ws.Range("I80").Value = 3478
$ws.Range("J80").Value = 3999.6667
$ws.Range("K80").Value = 3478
$ws.Range("L80").Value = 3999.6667
$ws.Range("M80").Value = -2480
$ws.Range("N80").Value = -5995.6667

# Row 83
$ws.Range("H83").Value = 3673.625
$ws.Range("I83").Value = 3478
$ws.Range("J83").Value = 3999.6667
$ws.Range("K83").Value = 17390
$ws.Range("L83").Value = 19998.3335
$ws.Range("M83").Value = -12398
$ws.Range("N83").Value = -29982.3335

# Row 102
$ws.Range("H102").Value = 1260
$ws.Range("I102").Value = 1125.5454
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 1125.5454
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = 496.4546
$ws.Range("N102").Value = -5243.5

# Row 126
$ws.Range("H126").Value = 4289.067
$ws.Range("I126").Value = 4641.4546
$ws.Range("K126").Value = 13924.3638
$ws.Range("M126").Value = -11454.3638

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 20
$ws.Range("H20").Value = 34974.75
$ws.Range("I20").Value = 36633
$ws.Range("J20").Value = 30000
$ws.Range("K20").Value = 36633
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = -36407
$ws.Range("N20").Value = -30452

# Row 46
$ws.Range("H46").Value = 1429.4
$ws.Range("I46").Value = 1429.4
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1429.4
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1241.4
$ws.Range("N46").ClearContents()

# Row 68
$ws.Range("H68").Value = 1424.375
$ws.Range("I68").Value = 1424.375
$ws.Range("K68").Value = 1424.375
$ws.Range("M68").Value = -675.375

# Row 71
$ws.Range("H71").Value = 1424.375
$ws.Range("I71").Value = 1424.375
$ws.Range("K71").Value = 7121.875
$ws.Range("M71").Value = -3377.875

# Row 82
$ws.Range("H82").Value = 2283.3333
$ws.Range("I82").Value = 1900
$ws.Range("J82").Value = 2666.6667
$ws.Range("K82").Value = 1900
$ws.Range("L82").Value = 2666.6667
$ws.Range("M82").Value = -1539
$ws.Range("N82").Value = -3388.6667

# Row 85
$ws.Range("H85").Value = 2283.3333
$ws.Range("I85").Value = 1900
$ws.Range("J85").Value = 2666.6667
$ws.Range("K85").Value = 1900
$ws.Range("L85").Value = 2666.6667
$ws.Range("M85").Value = -652
$ws.Range("N85").Value = -5162.6667

# Row 136
$ws.Range("H136").Value = 1873.5454
$ws.Range("I136").Value = 1623.2858
$ws.Range("J136").Value = 1990.3334
$ws.Range("K136").Value = 4869.857400000001
$ws.Range("L136").Value = 5971.0002
$ws.Range("M136").Value = -2319.857400000001
$ws.Range("N136").Value = -11071.0002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 122
$ws.Range("H122").Value = 346683.78
$ws.Range("I122").Value = 401883.1
$ws.Range("K122").Value = 1205649.3
$ws.Range("M122").Value = -1203199.3

# Row 136
$ws.Range("H136").Value = 11113255
$ws.Range("I136").Value = 11906834
$ws.Range("K136").Value = 35720502
$ws.Range("M136").Value = -35717952
